# The author removed the "「ツイッターでフォローして下さい」" entry (row 135)
# from the posts sheet. Deleting the entire row shifts every following
# row up by one, which also updates the sheet's used-range dimension
# from A1:C188 to A1:C187 automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(135).EntireRow.Delete()
